$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header: "age" -> "Age"
$ws.Range("B1").Value = "Age"

# New data for the LSM-analysis dataset "all_files_01_2024"
# Columns: A=RegistryCode, B=Age, C=WSO, D=WPM, E=FPM
$data = @(
    @("sub 1419", 66, 823.4285714285714, 0.7552083333333334, 0.6190476190476191),
    @("sub 3058", 61, 1422.428571428571, 0.7526881720430108, 0.6440677966101694),
    @("sub 2917", 59, 450.8571428571428, 0.9895833333333334, 0.8723404255319149),
    @("sub 3104", 50, 405.5714285714286, 0.7684210526315789, 0.7315789473684211),
    @("sub 3889", 62, 202.7142857142857, 0.824468085106383, 0.7819148936170213),
    @("sub 3642", 60, 275.0, 0.7894736842105263, 0.7604166666666666),
    @("sub 3035", 67, 843.2857142857143, 0.7777777777777778, 0.6770833333333334),
    @("sub 4182", 53, 146.2857142857143, 0.8952879581151832, 0.8663101604278075),
    @("sub 4137", 70, 162.2857142857143, 0.734375, 0.7120418848167539),
    @("sub 3583", 51, 481.5714285714286, 0.531578947368421, 0.4764397905759162),
    @("sub 4281", 46, 128.5714285714286, 0.774869109947644, 0.6613756613756614),
    @("sub 4191", 55, 131.4285714285714, 0.6596858638743456, 0.7225130890052356),
    @("sub 3201", 51, 436.4285714285714, 0.680628272251309, 0.6210526315789474),
    @("sub 4275", 55, 121.5714285714286, 0.5668449197860963, 0.6243386243386243),
    @("sub 4170", 58, 308.0, 0.8105263157894737, 0.6894736842105263),
    @("sub 3676", 63, 280.8571428571428, 0.8, 0.7671957671957672),
    @("sub 4431", 71, 96.71428571428571, 0.6701570680628273, 0.6117021276595744),
    @("sub 4294", 67, 107.5714285714286, 0.6666666666666666, 0.6031746031746031),
    @("sub 3650", 62, 331.5714285714286, 0.7513227513227513, 0.601063829787234),
    @("sub 4140", 44, 161.7142857142857, 0.7368421052631579, 0.6648351648351648),
    @("sub 3912", 55, 229.8571428571429, 0.6631578947368421, 0.6210526315789474),
    @("sub 4051", 59, 196.7142857142857, 0.8109756097560976, 0.632),
    @("sub 3286", 76, 390.1428571428572, 0.8691099476439791, 0.8148148148148148),
    @("sub 4208", 51, 137.2857142857143, 0.78125, 0.7120418848167539),
    @("sub 1536", 69, 1061.0, 0.6842105263157895, 0.6421052631578947),
    @("sub 4439", 56, 79.14285714285714, 0.746031746031746, 0.6455026455026455),
    @("sub 4466", 65, 95.85714285714286, 0.7552083333333334, 0.6858638743455497),
    @("sub 4504", 55, 78.14285714285714, 0.675531914893617, 0.6296296296296297),
    @("sub 3887", 68, 310.4285714285714, 0.9162303664921466, 0.7789473684210526),
    @("sub 4391", 61, 101.0, 0.7789473684210526, 0.6808510638297872),
    @("sub 4546", 40, 123.5714285714286, 0.9947916666666666, 0.9479166666666666),
    @("sub 4227", 38, 186.7142857142857, 0.1588785046728972, 0.3565217391304348),
    @("sub 4532", 63, 89.42857142857143, 0.8429319371727748, 0.7083333333333334),
    @("sub 3299", 60, 509.2857142857143, 0.9166666666666666, 0.8138297872340425),
    @("sub 4720", 59, 43.28571428571428, 0.7172774869109948, 0.6631578947368421),
    @("sub 4607", 63, 60.42857142857143, 0.6526315789473685, 0.5925925925925926),
    @("sub 4728", 43, 34.42857142857143, 0.837696335078534, 0.7277486910994765),
    @("sub 4769", 54, 31.71428571428572, 0.6797752808988764, 0.5689655172413793),
    @("sub 4777", 51, 28.71428571428572, 0.8105263157894737, 0.7421052631578947),
    @("sub 3396", 54, 337.2857142857143, 0.5789473684210527, 0.5578947368421052),
    @("sub 4467", 60, 31.28571428571428, 0.5585106382978723, 0.510752688172043),
    @("sub 3392", 53, 331.1428571428572, 0.6157894736842106, 0.6458333333333334),
    @("sub 3960", 79, 239.7142857142857, 0.4658385093167702, 0.5098039215686274),
    @("sub 3154", 30, 443.1428571428572, 0.796875, 0.75),
    @("sub 3994", 55, 179.8571428571429, 0.4789473684210526, 0.4789473684210526),
    @("sub 4086", 59, 141.7142857142857, 0.9739583333333334, 0.8631578947368421),
    @("sub 4002", 69, 185.2857142857143, 0.7368421052631579, 0.7105263157894737),
    @("sub 2522", 40, 555.8571428571429, 0.8324607329842932, 0.7947368421052632),
    @("sub 3364", 34, 334.4285714285714, 0.9162303664921466, 0.6684491978609626),
    @("sub 2884", 55, 516.8571428571429, 0.5842105263157895, 0.5882352941176471),
    @("sub 2998", 59, 469.2857142857143, 0.93048128342246, 0.848314606741573)
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $rowIndex++
}
